$d = $word.ActiveDocument

$replacements = @(
    @{old="166×4="; new="627×3="},
    @{old="188×3="; new="125×3="},
    @{old="972×3="; new="526×8="},
    @{old="236×8="; new="291×2="},
    @{old="727×8="; new="335×4="},
    @{old="949×6="; new="304×7="},
    @{old="142×5="; new="363×3="},
    @{old="705×9="; new="607×2="},
    @{old="187×7="; new="971×5="},
    @{old="420×7="; new="687×3="},
    @{old="751×3="; new="364×3="},
    @{old="462×6="; new="448×8="},
    @{old="907×4="; new="379×8="},
    @{old="222×3="; new="499×7="},
    @{old="707×4="; new="744×7="},
    @{old="273×5="; new="963×7="},
    @{old="205×9="; new="852×2="},
    @{old="976×5="; new="154×5="},
    @{old="540×6="; new="246×3="},
    @{old="268×3="; new="949×2="},
    @{old="390×9="; new="408×7="},
    @{old="524×6="; new="390×3="},
    @{old="220×7="; new="669×2="},
    @{old="847×8="; new="732×2="},
    @{old="729×5="; new="854×9="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}

Write-Host "Done applying $($replacements.Count) replacements"
